$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old SubjectAverage/SubjectGrades rows 12-13)
# so the table shrinks from 8 data/summary rows to 6, shifting nothing else.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# Row 6: Steve Coup (file no. 15613)
$ws.Range("A6").Value = 15613
$ws.Range("B6").Value = "Steve Coup"
$ws.Range("C6").Value = 80
$ws.Range("D6").Value = 44
$ws.Range("E6").Value = 68
$ws.Range("F6").Value = 48
$ws.Range("G6").Value = 70
$ws.Range("H6").Value = 81
$ws.Range("I6").Value = 54
$ws.Range("J6").Value = 89
$ws.Range("K6").Value = 534
$ws.Range("L6").Value = 66.75
$ws.Range("M6").Value = "B-"
$ws.Range("N6").Value = 1

# Row 7: Fredrick Ndote (file no. 15611) - name/file no unchanged, scores updated
$ws.Range("A7").Value = 15611
$ws.Range("B7").Value = "Fredrick Ndote"
$ws.Range("C7").Value = 70
$ws.Range("D7").Value = 64
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 59
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 76
$ws.Range("I7").Value = 36
$ws.Range("J7").Value = 86
$ws.Range("K7").Value = 531
$ws.Range("L7").Value = 66.375
$ws.Range("M7").Value = "B-"
$ws.Range("N7").Value = 2

# Row 8: Bostwald Kite (file no. 15612)
$ws.Range("A8").Value = 15612
$ws.Range("B8").Value = "Bostwald Kite"
$ws.Range("C8").Value = 56
$ws.Range("D8").Value = 79
$ws.Range("E8").Value = 36
$ws.Range("F8").Value = 60
$ws.Range("G8").Value = 87
$ws.Range("H8").Value = 48
$ws.Range("I8").Value = 68
$ws.Range("J8").Value = 69
$ws.Range("K8").Value = 503
$ws.Range("L8").Value = 62.875
$ws.Range("M8").Value = "C+"
$ws.Range("N8").Value = 3

# Row 9: SubjectTotal (no file no., no rank)
$ws.Range("A9").Value = "SubjectTotal"
$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = 206
$ws.Range("D9").Value = 187
$ws.Range("E9").Value = 154
$ws.Range("F9").Value = 167
$ws.Range("G9").Value = 247
$ws.Range("H9").Value = 205
$ws.Range("I9").Value = 158
$ws.Range("J9").Value = 244
$ws.Range("K9").Value = 1568
$ws.Range("L9").Value = 196
$ws.Range("M9").Value = "-"
$ws.Range("N9").Value = ""

# Row 10: SubjectAverage
$ws.Range("A10").Value = "SubjectAverage"
$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = 68.66666666666667
$ws.Range("D10").Value = 62.33333333333334
$ws.Range("E10").Value = 51.33333333333334
$ws.Range("F10").Value = 55.66666666666666
$ws.Range("G10").Value = 82.33333333333333
$ws.Range("H10").Value = 68.33333333333333
$ws.Range("I10").Value = 52.66666666666666
$ws.Range("J10").Value = 81.33333333333333
$ws.Range("K10").Value = 522.6666666666666
$ws.Range("L10").Value = 65.33333333333333
$ws.Range("M10").Value = "-"
$ws.Range("N10").Value = ""

# Row 11: SubjectGrades
$ws.Range("A11").Value = "SubjectGrades"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "B-"
$ws.Range("D11").Value = "C+"
$ws.Range("E11").Value = "C-"
$ws.Range("F11").Value = "C"
$ws.Range("G11").Value = "A-"
$ws.Range("H11").Value = "B-"
$ws.Range("I11").Value = "C-"
$ws.Range("J11").Value = "A-"
$ws.Range("K11").Value = "-"
$ws.Range("L11").Value = "B-"
$ws.Range("M11").Value = "-"
$ws.Range("N11").Value = ""
